$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCDTtiNTY")

# Update the three variables (LDVs, HDVs, motorbikes) plus the small
# correction to ships, as sent by Chris on 5/5/2020.
$ws.Range("B2").Value = 0.07692
$ws.Range("C2").Value = 0.085

$ws.Range("C3").Value = 0.045

$ws.Range("B6").Value = 0.0298
$ws.Range("C6").Value = 0.0298

$ws.Range("B7").Value = 0.11
$ws.Range("C7").Value = 0

# Make this sheet the active one and select B2:C7 with B2 as the active
# cell, matching the saved view state.
$ws.Activate()
$ws.Range("B2:C7").Select()
